$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage for numeric-looking Price values so Excel does not
# reinterpret strings like "1.003" as a number (matches original inlineStr text cells).
function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
}

# Row 2
Set-TextValue $ws.Range("D2") "28.481.35"
$ws.Range("E2").Value = "  -0.98%  "

# Row 3
Set-TextValue $ws.Range("D3") "1.818.92"
$ws.Range("E3").Value = "  -1.00%  "

# Row 4
$ws.Range("E4").Value = "  +0.67%  "

# Row 5
Set-TextValue $ws.Range("D5") "315.03"
$ws.Range("E5").Value = "  -0.89%  "

# Row 6
Set-TextValue $ws.Range("D6") "1.003"
$ws.Range("E6").Value = "  +0.46%  "

# Row 7
Set-TextValue $ws.Range("D7") "0.5053"
$ws.Range("E7").Value = "  -5.70%  "

# Row 8
Set-TextValue $ws.Range("D8") "0.3844"
$ws.Range("E8").Value = "  -3.22%  "

# Row 9
Set-TextValue $ws.Range("D9") "0.08500"
$ws.Range("E9").Value = "  +9.74%  "

# Row 10
$ws.Range("E10").Value = "  +0.23%  "

# Row 11
Set-TextValue $ws.Range("D11") "1.106"
$ws.Range("E11").Value = "  -1.85%  "

# Row 12
Set-TextValue $ws.Range("D12") "6.396"
$ws.Range("E12").Value = "  +0.18%  "

# Row 13
Set-TextValue $ws.Range("D13") "20.99"
$ws.Range("E13").Value = "  -1.27%  "

# Row 14
Set-TextValue $ws.Range("D14") "1.003"
$ws.Range("E14").Value = "  +0.36%  "

# Row 15
Set-TextValue $ws.Range("D15") "7.490"
$ws.Range("E15").Value = "  -1.64%  "

# Row 16
Set-TextValue $ws.Range("D16") "1.812.09"
$ws.Range("E16").Value = "  -1.11%  "

# Row 17
Set-TextValue $ws.Range("D17") "0.00001144"
$ws.Range("E17").Value = "  +4.82%  "

# Row 18
Set-TextValue $ws.Range("D18") "93.48"
$ws.Range("E18").Value = "  +1.90%  "

# Row 19
Set-TextValue $ws.Range("D19") "0.06665"
$ws.Range("E19").Value = "  +1.20%  "

# Row 20
Set-TextValue $ws.Range("D20") "17.69"
$ws.Range("E20").Value = "  -1.01%  "

# Row 21
$ws.Range("E21").Value = "  +0.23%  "

# Row 22
Set-TextValue $ws.Range("D22") "6.065"
$ws.Range("E22").Value = "  -0.68%  "

# Row 23
Set-TextValue $ws.Range("D23") "28.510.82"
$ws.Range("E23").Value = "  -0.77%  "

# Row 24
$ws.Range("E24").Value = "  +1.40%  "

# Row 25
Set-TextValue $ws.Range("D25") "2.269"
$ws.Range("E25").Value = "  +1.15%  "

# Row 26
Set-TextValue $ws.Range("D26") "21.20"
$ws.Range("E26").Value = "  +1.49%  "

# Row 27
Set-TextValue $ws.Range("D27") "157.88"
$ws.Range("E27").Value = "  +0.62%  "

# Row 28
Set-TextValue $ws.Range("D28") "2.023.58"
$ws.Range("E28").Value = "  -0.92%  "

# Row 29
Set-TextValue $ws.Range("D29") "2.383"
$ws.Range("E29").Value = "  -2.86%  "

# Row 30
Set-TextValue $ws.Range("D30") "126.05"
$ws.Range("E30").Value = "  +0.38%  "

# Row 31
Set-TextValue $ws.Range("D31") "1.106"
$ws.Range("E31").Value = "  -3.62%  "

# Row 32
Set-TextValue $ws.Range("D32") "0.1073"
$ws.Range("E32").Value = "  -4.15%  "

# Row 33
Set-TextValue $ws.Range("D33") "5.737"
$ws.Range("E33").Value = "  -0.75%  "

# Row 34
$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue $ws.Range("D34") "0.07496"
$ws.Range("E34").Value = "  +1.72%  "

# Row 35
$ws.Range("B35").Value = "HuobiToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextValue $ws.Range("D35") "3.678"
$ws.Range("E35").Value = "  +0.46%  "

# Row 36
Set-TextValue $ws.Range("D36") "0.2222"
$ws.Range("E36").Value = "  -2.17%  "

# Row 37
Set-TextValue $ws.Range("D37") "0.02355"
$ws.Range("E37").Value = "  -0.16%  "

# Row 38
Set-TextValue $ws.Range("D38") "5.212"
$ws.Range("E38").Value = "  -0.52%  "

# Row 39
Set-TextValue $ws.Range("D39") "8.696"
$ws.Range("E39").Value = "  -2.34%  "

# Row 40
Set-TextValue $ws.Range("D40") "0.6323"
$ws.Range("E40").Value = "  +0.03%  "

# Row 41
Set-TextValue $ws.Range("D41") "11.22"
$ws.Range("E41").Value = "  -2.28%  "

# Row 42
Set-TextValue $ws.Range("D42") "1.187"
$ws.Range("E42").Value = "  -0.83%  "

# Row 43
$ws.Range("E43").Value = "  +0.88%  "

# Row 44
Set-TextValue $ws.Range("D44") "13.62"
$ws.Range("E44").Value = "  +0.78%  "

# Row 45
$ws.Range("B45").Value = "PancakeSwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue $ws.Range("D45") "3.741"
$ws.Range("E45").Value = "  +0.69%  "

# Row 46
$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
Set-TextValue $ws.Range("D46") "0.5920"
$ws.Range("E46").Value = "  -0.26%  "

# Row 47
Set-TextValue $ws.Range("D47") "125.36"
$ws.Range("E47").Value = "  -0.43%  "

# Row 48
$ws.Range("E48").Value = "  -0.83%  "

# Row 49
Set-TextValue $ws.Range("D49") "1.198"
$ws.Range("E49").Value = "  -0.26%  "

# Row 50
Set-TextValue $ws.Range("D50") "0.06997"
$ws.Range("E50").Value = "  +0.53%  "

# Row 51
$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue $ws.Range("D51") "74.04"
$ws.Range("E51").Value = "  -0.75%  "
